$d = $word.ActiveDocument

# Remove the BARCODE field entirely (begin/instrText/separate/result/end runs),
# clearing any inherited character formatting along with it.
$f = $d.Fields.Item(1)
$f.Delete()

# Build the "Best practice" x11 text without relying on string repetition operators.
$bestPractice = ""
for ($i = 0; $i -lt 11; $i++) {
    $bestPractice = $bestPractice + "Best practice"
}

# Replace the (now empty) document content with both paragraphs in one shot so
# neither run picks up stray inherited run formatting.
$full = "Evaluation Only. Created with Aspose.Words. Copyright 2003-2016 Aspose Pty Ltd.`r" + $bestPractice
$d.Content.Text = $full

# Apply Bold / red color / 12pt to the first paragraph's text only (not its
# paragraph mark), matching the evaluation-watermark run formatting.
$p1 = $d.Paragraphs.Item(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$r1.Font.Bold = 1
$r1.Font.Color = 255
$r1.Font.Size = 12

# Update the page margins (values are in points; originals were all 1in/72pt).
$d.PageSetup.TopMargin = 56.7
$d.PageSetup.RightMargin = 42.5
$d.PageSetup.BottomMargin = 56.7
$d.PageSetup.LeftMargin = 85.05
